$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update 想去人数 (interest count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1341
$ws1.Range("F4").Value = 8

# Sheet "全部类型" (All types) - same updates mirrored
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1341
$ws4.Range("F4").Value = 8
